# Apply updated performance metrics to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - D-RIBose
$ws.Range("E2").Value = "225/509"
$ws.Range("F2").Value = 0.7822222222222223
$ws.Range("G2").Value = 0.704
$ws.Range("H2").Value = 0.773
$ws.Range("I2").Value = 0.741
$ws.Range("J2").Value = 0.758
$ws.Range("K2").Value = 0.757
$ws.Range("L2").Value = 0.7574054054054054

# Row 3 - N-AcetylGlucosamine
$ws.Range("E3").Value = "307/509"
$ws.Range("F3").Value = 0.758957654723127
$ws.Range("G3").Value = 0.8204225352112676
$ws.Range("H3").Value = 0.707
$ws.Range("I3").Value = 0.788
$ws.Range("J3").Value = 0.754
$ws.Range("K3").Value = 0.748
$ws.Range("L3").Value = 0.7457668231611894

# Row 4 - SALicin
$ws.Range("E4").Value = "201/509"
$ws.Range("F4").Value = 0.736318407960199
$ws.Range("G4").Value = 0.6271186440677966
$ws.Range("H4").Value = 0.757
$ws.Range("I4").Value = 0.677
$ws.Range("J4").Value = 0.723
$ws.Range("K4").Value = 0.717
$ws.Range("L4").Value = 0.7164897249643012

# Row 5 - D-CELlobiose (E5 unchanged)
$ws.Range("F5").Value = 0.7962962962962963
$ws.Range("G5").Value = 0.589041095890411
$ws.Range("H5").Value = 0.8070000000000001
$ws.Range("I5").Value = 0.677
$ws.Range("J5").Value = 0.758
$ws.Range("K5").Value = 0.742
$ws.Range("L5").Value = 0.7376239962210676

# Row 6 - D-LACtose (bovine origin)
$ws.Range("E6").Value = "165/509"
$ws.Range("F6").Value = 0.8181818181818182
$ws.Range("G6").Value = 0.6367924528301887
$ws.Range("H6").Value = 0.833
$ws.Range("I6").Value = 0.716
$ws.Range("J6").Value = 0.79
$ws.Range("K6").Value = 0.775
$ws.Range("L6").Value = 0.7678911759100439

# Row 7 - D-MELibiose
$ws.Range("E7").Value = "176/509"
$ws.Range("F7").Value = 0.8068181818181818
$ws.Range("G7").Value = 0.6926829268292682
$ws.Range("H7").Value = 0.848
$ws.Range("I7").Value = 0.745
$ws.Range("J7").Value = 0.8090000000000001
$ws.Range("K7").Value = 0.797
$ws.Range("L7").Value = 0.7904204107830552

# Row 8 - D-SACcharose (sucrose)
$ws.Range("E8").Value = "325/509"
$ws.Range("F8").Value = 0.7569230769230769
$ws.Range("G8").Value = 0.8145695364238411
$ws.Range("H8").Value = 0.655
$ws.Range("I8").Value = 0.785
$ws.Range("J8").Value = 0.735
$ws.Range("K8").Value = 0.72
$ws.Range("L8").Value = 0.7164635121732732

# Row 9 - D-TREhalose
$ws.Range("E9").Value = "194/509"
$ws.Range("F9").Value = 0.8195876288659794
$ws.Range("G9").Value = 0.6943231441048034
$ws.Range("H9").Value = 0.824
$ws.Range("I9").Value = 0.752
$ws.Range("J9").Value = 0.794
$ws.Range("K9").Value = 0.788
$ws.Range("L9").Value = 0.7846615720524017
